$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set the Runmode column (C) to "Y" for all test rows so that all test
# cases are run.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update the active selection to match the new state.
$ws.Activate()
$ws.Range("C2:C7").Select()
